$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("H1").Value = "Hint"
$ws.Range("I1").Value = "Popup"

# Give them the same look as the rest of the header row (blue fill / white text),
# but without the border/wrap so they render like a plain header cell.
$headerCell = $ws.Range("G1")
$newHeaders = $ws.Range("H1:I1")
$newHeaders.Style = $headerCell.Style

# Move the active selection, matching where the editor ended up
$ws.Range("I3").Select()
